$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "29.100.49"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.11%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.837.99"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.08%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.9969"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.27%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "243.17"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.41%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.6244"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.26%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.27%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.07506"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.00%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.2946"
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "23.34"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +2.11%  "

$ws.Cells.Item(11, 5).Value = "  -0.67%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "1.837.04"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.01%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.020"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.43%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "0.6764"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.77%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "83.04"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.25%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "0.000009376"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -4.58%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "5.980"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -2.31%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "29.099.45"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.03%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "2.080.61"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.29%  "

$ws.Cells.Item(20, 5).Value = "  +0.99%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "226.75"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.18%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "0.9992"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.15%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "7.164"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -1.18%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "0.9981"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.28%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "160.10"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.31%  "

$ws.Cells.Item(26, 5).Value = "  -0.45%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "8.541"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.11%  "

$ws.Cells.Item(28, 5).Value = "  -0.38%  "

$ws.Cells.Item(29, 5).Value = "  -0.59%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "4.183"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +1.39%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "4.149"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +2.31%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.05571"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +3.84%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.205"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.30%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.7488"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.29%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.847"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.84%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "1.146"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.23%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "2.661"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.36%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "1.237.31"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.05%  "

$ws.Cells.Item(39, 5).Value = "  +0.35%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.01785"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.76%  "

$ws.Cells.Item(41, 5).Value = "  -0.88%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.9004"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.59%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.9986"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.32%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "102.44"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.20%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "1.984.88"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.09%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "66.53"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +2.45%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.00000000123"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.74%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.5076"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.78%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.4076"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.51%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "9.086"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.05%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.05839"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.57%  "
